$wb = $excel.ActiveWorkbook

$sheetsData = @(
    @{
        Name = "ROW35-FE-LIFTER"
        A = "2025-03-06 08:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Name = "ROW35-MID-LIFTER"
        A = "2025-03-06 08:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Name = "ROW02-FE-LIFTER"
        A = "2025-03-06 08:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Name = "ROW02-MID-LIFTER"
        A = "2025-03-06 08:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($row in $sheetsData) {
    $ws = $wb.Worksheets.Item($row.Name)
    $ws.Range("A49").Value = $row.A
    $ws.Range("B49").Value = $row.B
    $ws.Range("C49").Value = $row.C
    $ws.Range("D49").Value = $row.D
    $ws.Range("E49").Value = $row.E
    $ws.Range("F49").Value = $row.F

    # G column holds a number too large to round-trip through a double
    # without losing precision, so the source file stores it as text.
    # Force text formatting before assignment so COM doesn't coerce it
    # back into a floating point number.
    $ws.Range("G49").NumberFormat = "@"
    $ws.Range("G49").Value = $row.G

    $ws.Range("H49").Value = $row.H
    $ws.Range("I49").Value = $row.I
}
